# Update column F (dSF) values to reflect repulled data / recalculated mean values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 11
$ws.Range("F3").Value = 0
$ws.Range("F4").Value = -5
$ws.Range("F5").Value = 4
$ws.Range("F7").Value = -1
$ws.Range("F8").Value = -4
$ws.Range("F12").Value = -2
$ws.Range("F15").Value = -7
